$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("A9").Value = "KIRAN KUMAR"
$ws.Range("B9").Value = "OS"
$ws.Range("C9").Value = "Ftth OS_01.12.2025.xlsx"
$ws.Range("D9").Value = "2025-12-02 12:39"
$ws.Range("E9").Value = "2025-12"
